$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09355383683086951
$ws.Range("H2").Value = -1.232090647012569
$ws.Range("I2").Value = -39.30609181603342
$ws.Range("G3").Value = 0.09315875159826537
$ws.Range("H3").Value = 40.94907802408473
$ws.Range("G4").Value = -0.03485944612852575
$ws.Range("H4").Value = -471.8234803106694
$ws.Range("G5").Value = -0.01639358031618353
$ws.Range("H5").Value = -22.43323429810065
$ws.Range("G6").Value = -0.2341564340652409
$ws.Range("H6").Value = -5.859402673931979
$ws.Range("G7").Value = -0.2104800223316823
$ws.Range("H7").Value = 15.7661080613405
$ws.Range("G8").Value = -0.3637366742125455
$ws.Range("H8").Value = 1.760727126773091
$ws.Range("G9").Value = -0.4057254922218577
$ws.Range("H9").Value = -1.778853779324516
$ws.Range("G10").Value = -0.007219825182621979
$ws.Range("H10").Value = -144.5533783605492
$ws.Range("G11").Value = 0.1005468291799932
$ws.Range("H11").Value = 725.8412765061602
$ws.Range("G12").Value = 0.2228508434658396
$ws.Range("H12").Value = -1.909406441983653
$ws.Range("G13").Value = 0.2632002203600229
$ws.Range("H13").Value = -0.05330155153606771
$ws.Range("G14").Value = -0.06773546902957142
$ws.Range("H14").Value = -607.9914343616121
$ws.Range("G15").Value = -0.005463744584387026
$ws.Range("H15").Value = -127.0660385255932
$ws.Range("G16").Value = 0.1418167598673289
$ws.Range("H16").Value = 20.18460599830405
$ws.Range("G17").Value = 0.1737968055846164
$ws.Range("H17").Value = -20.58013773337449
$ws.Range("G18").Value = 0.04277160040352979
$ws.Range("H18").Value = -29.26409299965166
$ws.Range("G19").Value = 0.09722995818674758
$ws.Range("H19").Value = 7.926195178855125
$ws.Range("G20").Value = -0.1660079485391075
$ws.Range("H20").Value = -14.07052595502557
$ws.Range("G21").Value = -0.1820990095074113
$ws.Range("H21").Value = 8.872988891963239
$ws.Range("G22").Value = 0.04966537345960819
$ws.Range("H22").Value = -8.683621488165691
$ws.Range("G23").Value = 0.05907778408778319
$ws.Range("H23").Value = 44.65615587360412
$ws.Range("G24").Value = 0.1450338375382942
$ws.Range("H24").Value = 25.3161842722803
$ws.Range("G25").Value = 0.1424787478059539
$ws.Range("H25").Value = -6.31084389681687
$ws.Range("G26").Value = 0.00691776711528883
$ws.Range("H26").Value = -86.91719460275957
$ws.Range("G27").Value = 0.02558299430154447
$ws.Range("H27").Value = -49.30728615956224
$ws.Range("G28").Value = 0.1596815804839949
$ws.Range("H28").Value = 4.427742161918706
$ws.Range("G29").Value = 0.184984367196109
$ws.Range("H29").Value = 8.366151681444043
$ws.Range("G30").Value = 0.01876653806405204
$ws.Range("H30").Value = -4.086401214161048
$ws.Range("G31").Value = 0.02898067054104423
$ws.Range("H31").Value = 198.6167782357273
$ws.Range("G32").Value = 0.02271605651033349
$ws.Range("H32").Value = -39.08805768298264
$ws.Range("G33").Value = -0.004345234279451735
$ws.Range("H33").Value = -116.6460863976168
$ws.Range("G34").Value = 0.1055534554179075
$ws.Range("H34").Value = -17.51574745224605
$ws.Range("G35").Value = 0.1712226295672531
$ws.Range("H35").Value = 33.08086454868662
$ws.Range("G36").Value = -0.01190553006132855
$ws.Range("H36").Value = -179.2018307082073
$ws.Range("G37").Value = 0.0003132051888267698
$ws.Range("H37").Value = -97.95484393217158
$ws.Range("G38").Value = -0.05457278894645743
$ws.Range("H38").Value = -2571.975414760152
$ws.Range("G39").Value = -0.02434570647430555
$ws.Range("H39").Value = 27.12876261625307
$ws.Range("G40").Value = 0.1575386575255652
$ws.Range("H40").Value = 6.770405757709627
$ws.Range("G41").Value = 0.1546372211322629
$ws.Range("H41").Value = -4.189602376706962
$ws.Range("G42").Value = 0.05597683064991897
$ws.Range("H42").Value = -13.30145789624524
$ws.Range("G43").Value = 0.01972873107971296
$ws.Range("H43").Value = -43.24386036991727
$ws.Range("G44").Value = 0.02978259451089352
$ws.Range("H44").Value = 111.0327580943088
$ws.Range("G45").Value = 0.03904622824837738
$ws.Range("H45").Value = -4.899918042726817
$ws.Range("G46").Value = -0.05041354899459961
$ws.Range("H46").Value = 23.40678646048049
$ws.Range("G47").Value = -0.08321886499833031
$ws.Range("H47").Value = -101.4511345678969
$ws.Range("G48").Value = -0.131274698122713
$ws.Range("H48").Value = -4.206605657500666
$ws.Range("G49").Value = -0.1446512613410191
$ws.Range("H49").Value = 26.75195916328642
$ws.Range("G50").Value = 0.08833056891999445
$ws.Range("H50").Value = -18.86873387165947
$ws.Range("G51").Value = 0.1331291343494277
$ws.Range("H51").Value = 32.76955677321372
$ws.Range("G52").Value = 0.03351564910520254
$ws.Range("H52").Value = -43.78366098153749
$ws.Range("G53").Value = 0.05987009946767346
$ws.Range("H53").Value = -11.35960392184528
$ws.Range("G54").Value = -0.03872326659956191
$ws.Range("H54").Value = 44.61821338309785
$ws.Range("G55").Value = -0.07296990268570148
$ws.Range("H55").Value = 5.523704875532832
$ws.Range("G56").Value = 0.1163191734599165
$ws.Range("H56").Value = 153.8274923727075
$ws.Range("G57").Value = 0.06219180581000838
$ws.Range("H57").Value = 1102.871645863639
